# Auto-generated edit script: updates cached numeric values on the
# "Omega_Profits" data sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the target snapshot. Columns H-N hold plain numbers (no
# formulas), so each touched cell is just re-written with its new value;
# a couple of cells are newly populated (previously empty) or cleared
# out entirely (previously populated), which ClearContents() mirrors.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1474
$ws.Range("I6").Value = 430
$ws.Range("K6").Value = 1290
$ws.Range("M6").Value = -1178
# Row 17
$ws.Range("H17").Value = 3325.389
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 3678.5625
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 11035.6875
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -11371.6875
# Row 43
$ws.Range("H43").Value = 2406.7144
$ws.Range("J43").Value = 2474.25
$ws.Range("L43").Value = 2474.25
$ws.Range("N43").Value = -2612.25
# Row 74
$ws.Range("H74").Value = 3453.3845
$ws.Range("I74").Value = 3453.3845
$ws.Range("K74").Value = 3453.3845
$ws.Range("M74").Value = -2517.3845
# Row 77
$ws.Range("H77").Value = 3453.3845
$ws.Range("I77").Value = 3453.3845
$ws.Range("K77").Value = 17266.9225
$ws.Range("M77").Value = -12586.9225
# Row 80
$ws.Range("H80").Value = 806.9231
$ws.Range("I80").Value = 399.8
$ws.Range("J80").Value = 1061.375
$ws.Range("K80").Value = 1199.4
$ws.Range("L80").Value = 3184.125
$ws.Range("M80").Value = -201.4000000000001
$ws.Range("N80").Value = -5180.125
# Row 83
$ws.Range("H83").Value = 806.9231
$ws.Range("I83").Value = 399.8
$ws.Range("J83").Value = 1061.375
$ws.Range("K83").Value = 3598.2
$ws.Range("L83").Value = 9552.375
$ws.Range("M83").Value = 1393.8
$ws.Range("N83").Value = -19536.375
# Row 106
$ws.Range("H106").Value = 7967.7856
$ws.Range("I106").Value = 7967.7856
$ws.Range("K106").Value = 7967.7856
$ws.Range("M106").Value = -7336.7856
# Row 132
$ws.Range("H132").Value = 2572.4753
$ws.Range("I132").Value = 2261.7896
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 6785.3688
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -4255.3688
$ws.Range("N132").Value = -26059.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6961.8184
$ws.Range("I32").Value = 1731.125
$ws.Range("K32").Value = 1731.125
$ws.Range("M32").Value = -1444.125
# Row 61
$ws.Range("H61").Value = 3911.5435
$ws.Range("I61").Value = 3148.925
$ws.Range("J61").Value = 8995.666999999999
$ws.Range("K61").Value = 3148.925
$ws.Range("L61").Value = 8995.666999999999
$ws.Range("M61").Value = -2936.925
$ws.Range("N61").Value = -9419.666999999999
# Row 74
$ws.Range("H74").Value = 4133
$ws.Range("I74").Value = 3586.8462
$ws.Range("J74").Value = 6499.6665
$ws.Range("K74").Value = 3586.8462
$ws.Range("L74").Value = 6499.6665
$ws.Range("M74").Value = -2712.8462
$ws.Range("N74").Value = -8247.666499999999
# Row 77
$ws.Range("H77").Value = 4133
$ws.Range("I77").Value = 3586.8462
$ws.Range("J77").Value = 6499.6665
$ws.Range("K77").Value = 17934.231
$ws.Range("L77").Value = 32498.3325
$ws.Range("M77").Value = -13566.231
$ws.Range("N77").Value = -41234.3325
# Row 136
$ws.Range("H136").Value = 3911.5435
$ws.Range("I136").Value = 3148.925
$ws.Range("J136").Value = 8995.666999999999
$ws.Range("K136").Value = 9446.775000000001
$ws.Range("L136").Value = 26987.001
$ws.Range("M136").Value = -6896.775000000001
$ws.Range("N136").Value = -32087.001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 45000
$ws.Range("J60").Value = 45000
$ws.Range("L60").Value = 45000
$ws.Range("N60").Value = -46198
# Row 80
$ws.Range("H80").Value = 715.5294
$ws.Range("I80").Value = 132.75
$ws.Range("K80").Value = 132.75
$ws.Range("M80").Value = 865.25
# Row 83
$ws.Range("H83").Value = 715.5294
$ws.Range("I83").Value = 132.75
$ws.Range("K83").Value = 663.75
$ws.Range("M83").Value = 4328.25
# Row 86
$ws.Range("H86").Value = 4389059.5
$ws.Range("I86").Value = 7578648
$ws.Range("K86").Value = 7578648
$ws.Range("M86").Value = -7577525
# Row 89
$ws.Range("H89").Value = 4389059.5
$ws.Range("I89").Value = 7578648
$ws.Range("K89").Value = 37893240
$ws.Range("M89").Value = -37887624
# Row 94
$ws.Range("H94").Value = 26804506
$ws.Range("I94").Value = 32609448
$ws.Range("J94").Value = 101774.4
$ws.Range("K94").Value = 32609448
$ws.Range("L94").Value = 101774.4
$ws.Range("M94").Value = -32608997
$ws.Range("N94").Value = -102676.4
# Row 105
$ws.Range("H105").Value = 3550.4
$ws.Range("I105").Value = 2853.3333
$ws.Range("K105").Value = 2853.3333
$ws.Range("M105").Value = -1106.3333
# Row 107
$ws.Range("H107").Value = 2933.5715
$ws.Range("I107").Value = 2905.5
$ws.Range("J107").Value = 3036.5
$ws.Range("K107").Value = 2905.5
$ws.Range("L107").Value = 3036.5
$ws.Range("M107").Value = -985.5
$ws.Range("N107").Value = -6876.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 698.3077
$ws.Range("I22").Value = 698.3077
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 698.3077
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -348.3077
$ws.Range("N22").ClearContents()
# Row 31
$ws.Range("H31").Value = 5741.745
$ws.Range("I31").Value = 5422.5713
$ws.Range("K31").Value = 5422.5713
$ws.Range("M31").Value = -5127.5713
# Row 34
$ws.Range("H34").Value = 5741.745
$ws.Range("I34").Value = 5422.5713
$ws.Range("K34").Value = 5422.5713
$ws.Range("M34").Value = -5220.5713
# Row 94
$ws.Range("H94").Value = 1504.3334
$ws.Range("I94").Value = 1549.5
$ws.Range("K94").Value = 1549.5
$ws.Range("M94").Value = -1098.5
# Row 122
$ws.Range("H122").Value = 86153.55499999999
$ws.Range("I122").Value = 106328.484
$ws.Range("J122").Value = 2571.7144
$ws.Range("K122").Value = 318985.452
$ws.Range("L122").Value = 7715.1432
$ws.Range("M122").Value = -316535.452
$ws.Range("N122").Value = -12615.1432
# Row 131
$ws.Range("H131").Value = 46995.2
$ws.Range("J131").Value = 47773.223
$ws.Range("L131").Value = 47773.223
$ws.Range("N131").Value = -57853.223
# Row 132
$ws.Range("H132").Value = 1924.2
$ws.Range("I132").Value = 1648.5
$ws.Range("J132").Value = 2475.6
$ws.Range("K132").Value = 4945.5
$ws.Range("L132").Value = 7426.799999999999
$ws.Range("M132").Value = -2415.5
$ws.Range("N132").Value = -12486.8
# Row 134
$ws.Range("H134").Value = 1154.4
$ws.Range("I134").Value = 1069.25
$ws.Range("K134").Value = 3207.75
$ws.Range("M134").Value = -672.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 496.25
$ws.Range("I5").Value = 381.7143
$ws.Range("J5").Value = 1298
$ws.Range("K5").Value = 1145.1429
$ws.Range("L5").Value = 3894
$ws.Range("M5").Value = -1033.1429
$ws.Range("N5").Value = -4118
# Row 135
$ws.Range("H135").Value = 496.25
$ws.Range("I135").Value = 381.7143
$ws.Range("J135").Value = 1298
$ws.Range("K135").Value = 3435.4287
$ws.Range("L135").Value = 11682
$ws.Range("M135").Value = -900.4286999999999
$ws.Range("N135").Value = -16752

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5974.5864
$ws.Range("I80").Value = 5346.077
$ws.Range("J80").Value = 6485.25
$ws.Range("K80").Value = 5346.077
$ws.Range("L80").Value = 6485.25
$ws.Range("M80").Value = -4348.077
$ws.Range("N80").Value = -8481.25
# Row 83
$ws.Range("H83").Value = 5974.5864
$ws.Range("I83").Value = 5346.077
$ws.Range("J83").Value = 6485.25
$ws.Range("K83").Value = 26730.385
$ws.Range("L83").Value = 32426.25
$ws.Range("M83").Value = -21738.385
$ws.Range("N83").Value = -42410.25
# Row 99
$ws.Range("H99").Value = 25374.5
$ws.Range("J99").Value = 69999
$ws.Range("L99").Value = 69999
$ws.Range("N99").Value = -74491
# Row 126
$ws.Range("H126").Value = 6021.069
$ws.Range("I126").Value = 5382.35
$ws.Range("J126").Value = 7440.4443
$ws.Range("K126").Value = 16147.05
$ws.Range("L126").Value = 22321.3329
$ws.Range("M126").Value = -13677.05
$ws.Range("N126").Value = -27261.3329

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3033.3333
$ws.Range("I46").Value = 2333.3333
$ws.Range("J46").Value = 3733.3333
$ws.Range("K46").Value = 2333.3333
$ws.Range("L46").Value = 3733.3333
$ws.Range("M46").Value = -2145.3333
$ws.Range("N46").Value = -4109.3333
# Row 130
$ws.Range("H130").Value = 92625.14
$ws.Range("J130").Value = 92625.14
$ws.Range("L130").Value = 92625.14
$ws.Range("N130").Value = -102665.14
# Row 136
$ws.Range("H136").Value = 6232
$ws.Range("I136").Value = 6938.4
$ws.Range("K136").Value = 20815.2
$ws.Range("M136").Value = -18265.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 14999.5
$ws.Range("I81").Value = 14999.5
$ws.Range("K81").Value = 29999
$ws.Range("M81").Value = -28938
# Row 84
$ws.Range("H84").Value = 14999.5
$ws.Range("I84").Value = 14999.5
$ws.Range("K84").Value = 149995
$ws.Range("M84").Value = -144691
# Row 132
$ws.Range("H132").Value = 5625.8687
$ws.Range("I132").Value = 3448.121
$ws.Range("K132").Value = 10344.363
$ws.Range("M132").Value = -7814.363000000001
# Row 133
$ws.Range("H133").Value = 91000
$ws.Range("J133").Value = 91000
$ws.Range("L133").Value = 91000
$ws.Range("N133").Value = -101120
# Row 136
$ws.Range("H136").Value = 4224.25
$ws.Range("I136").Value = 4311.148
$ws.Range("J136").Value = 3755
$ws.Range("K136").Value = 12933.444
$ws.Range("L136").Value = 11265
$ws.Range("M136").Value = -10383.444
$ws.Range("N136").Value = -16365

Write-Output "Done updating Omega_Profits sheets"
